$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 ("2021年") mirrors the layout of row 11 ("2020年").
# Copy row 11's column-A cell first so the new label cell inherits the
# same style (bold / centered / bordered) used by every other year label,
# then overwrite the copied value with the real 2021 label.
$ws.Cells.Item(11, 1).Copy($ws.Cells.Item(12, 1))
$ws.Cells.Item(12, 1).Value = "2021年"

$ws.Cells.Item(12, 4).Value = 7412
$ws.Cells.Item(12, 5).Value = 113344
$ws.Cells.Item(12, 6).Value = 41897
$ws.Cells.Item(12, 7).Value = 623
$ws.Cells.Item(12, 10).Value = 292328
$ws.Cells.Item(12, 11).Value = 6778
$ws.Cells.Item(12, 12).Value = 584
$ws.Cells.Item(12, 13).Value = 5
$ws.Cells.Item(12, 14).Value = 234
$ws.Cells.Item(12, 15).Value = 52
$ws.Cells.Item(12, 16).Value = 2692
$ws.Cells.Item(12, 18).Value = 120
$ws.Cells.Item(12, 19).Value = 75
$ws.Cells.Item(12, 20).Value = 5
$ws.Cells.Item(12, 21).Value = 35
$ws.Cells.Item(12, 22).Value = 182
$ws.Cells.Item(12, 23).Value = 251
$ws.Cells.Item(12, 24).Value = 458
$ws.Cells.Item(12, 25).Value = 786
$ws.Cells.Item(12, 27).Value = 1321
$ws.Cells.Item(12, 29).Value = 105
$ws.Cells.Item(12, 30).Value = 352
$ws.Cells.Item(12, 31).Value = 41129
$ws.Cells.Item(12, 32).Value = 58988
$ws.Cells.Item(12, 34).Value = 276
$ws.Cells.Item(12, 35).Value = 1880
$ws.Cells.Item(12, 36).Value = 24
$ws.Cells.Item(12, 37).Value = 230
$ws.Cells.Item(12, 39).Value = 186
$ws.Cells.Item(12, 40).Value = 777
$ws.Cells.Item(12, 41).Value = 534
$ws.Cells.Item(12, 42).Value = 250
$ws.Cells.Item(12, 43).Value = 1022
$ws.Cells.Item(12, 44).Value = 3
$ws.Cells.Item(12, 45).Value = 186
$ws.Cells.Item(12, 46).Value = 2037
$ws.Cells.Item(12, 47).Value = 5148
$ws.Cells.Item(12, 48).Value = 45
$ws.Cells.Item(12, 49).Value = 2273
$ws.Cells.Item(12, 50).Value = 29
